# A new weekly data point for Brocoli (Macroferia Regional de Talca) is inserted
# as a new row right before the current row 149. This pushes the existing rows
# 149..258 down by one (to 150..259), and the brand new row 149 gets its own
# Fecha (D) and Volumen (J) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 149; rows 149-258 shift down to 150-259.
$ws.Rows.Item(149).Insert()

# Populate the new row 149 with the data for the newly added record. All of the
# columns other than Fecha (D) and Volumen (J) match what used to be in row 149
# before the shift (Calidad/Precios/Origen, etc. are unchanged for this point).
$ws.Range("A149").Value = 5
$ws.Range("B149").Value = "Macroferia Regional de Talca"
$ws.Range("C149").Value = "Maule"
$ws.Range("D149").Value = 44574
$ws.Range("E149").Value = 7
$ws.Range("F149").Value = 100112023
$ws.Range("G149").Value = "Brócoli"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 3000
$ws.Range("K149").Value = 500
$ws.Range("L149").Value = 500
$ws.Range("M149").Value = 500
$ws.Range("N149").Value = "`$/unidad"
$ws.Range("O149").Value = "Región del Maule"
$ws.Range("P149").Value = 500
$ws.Range("Q149").Value = 1
$ws.Range("R149").Value = "Hortaliza"

Write-Host "Row inserted and populated"
